$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Fuel Filter" service row (row 21) that had been left blank.
# Order matches the original author's edit so new shared strings line up.
$ws.Range("B21").Value = "Fuel Filter"
$ws.Range("C21").Value = "kecil"
$ws.Range("D21").Value = "3 pcs"
$ws.Range("F21").Value = "25/3/2024"
$ws.Range("G21").Value = "25/6/2024"
$ws.Range("E21").Value = "-"

# Update the active selection to match the edited workbook.
$ws.Range("E22").Select()
